$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws1.Range("F2").Value = 7599
$ws1.Range("F4").Value = 7787
$ws1.Range("F8").Value = 6488
$ws1.Range("F9").Value = 3335
$ws1.Range("F13").Value = 35
$ws1.Range("F14").Value = 36
$ws1.Range("F15").Value = 55
$ws1.Range("F20").Value = 315
$ws1.Range("F21").Value = 3771
$ws1.Range("F23").Value = 363
$ws1.Range("F25").Value = 278
$ws1.Range("F26").Value = 1429
$ws1.Range("F30").Value = 1730
$ws1.Range("F32").Value = 39
$ws1.Range("F33").Value = 52
$ws1.Range("F34").Value = 3562
$ws1.Range("F35").Value = 272
$ws1.Range("F36").Value = 272
$ws1.Range("F40").Value = 1379
$ws1.Range("F43").Value = 625

$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws2.Range("F8").Value = 56
$ws2.Range("F13").Value = 84

$ws4 = $wb.Worksheets.Item(4)  # 全部类型
$ws4.Range("F6").Value = 7599
$ws4.Range("F7").Value = 7787
$ws4.Range("F11").Value = 6488
$ws4.Range("F12").Value = 3335
$ws4.Range("F14").Value = 35
$ws4.Range("F15").Value = 55
$ws4.Range("F21").Value = 315
$ws4.Range("F22").Value = 3771
$ws4.Range("F26").Value = 363
$ws4.Range("F28").Value = 278
$ws4.Range("F29").Value = 1429
$ws4.Range("F33").Value = 1730
$ws4.Range("F35").Value = 39
$ws4.Range("F36").Value = 52
$ws4.Range("F37").Value = 84
$ws4.Range("F38").Value = 3562
$ws4.Range("F39").Value = 272
$ws4.Range("F40").Value = 272
$ws4.Range("F45").Value = 1379
$ws4.Range("F49").Value = 625
